$d = $word.ActiveDocument

# Locate the paragraph that ends with "ALTER TABLE reservadequartosleonardo.pessoa ADD senha varchar(100) NOT NULL;"
# which is the last SQL statement paragraph before the trailing bookmark/blank paragraphs.
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "ALTER TABLE reservadequartosleonardo.pessoa ADD senha varchar(100) NOT NULL;") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$insertPoint = $word.ActiveDocument.Range($target.Range.End, $target.Range.End)

$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="2F2F2F"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="2F2F2F"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>ALTER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>TABLE</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>reservadequartosleonardo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>pessoa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>CHANGE</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>idEndereco</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>id_endereco</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="C1AA6C"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>NOT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>NULL</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="EECC64"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="2F2F2F"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>ALTER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>TABLE</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>reservadequartosleonardo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>pessoa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>CHANGE</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>estadoCivil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>id_estado_civil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="C1AA6C"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>INT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>NOT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>NULL</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="EECC64"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="2F2F2F"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>ALTER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>TABLE</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>reservadequartosleonardo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>pessoa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>MODIFY</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>COLUMN</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9E9E9E"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>id_estado_civil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="C1AA6C"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>INT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>NOT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="AAAAAA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="739ECA"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>NULL</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="EECC64"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>;</w:t></w:r></w:p>'

$insertPoint.InsertXML($xmlFrag)

Write-Output "Inserted. New paragraph count:"
Write-Output $d.Paragraphs.Count
